# Update latest output (run 95)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 751.1665979999999
$schedule.Range("F2").Value = 12.42008263888889

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B7").Value = 35.87996

$detailed.Range("B8").Value = 57.50228

$detailed.Range("B9").Value = 58.03822
$detailed.Range("C9").Value = "historical"

$detailed.Range("B10").Value = 57.52849
$detailed.Range("C10").Value = "historical"

$detailed.Range("B11").Value = 61.43135

$detailed.Range("B12").Value = 60.04433

$detailed.Range("B13").Value = 65

$detailed.Range("B17").Value = 10.07303

$detailed.Range("B18").Value = 7.0734

$detailed.Range("B19").Value = 7.08291

$detailed.Range("B20").Value = 0.68371

$detailed.Range("B21").Value = -3.6481

$detailed.Range("B23").Value = -5.74313

$detailed.Range("B24").Value = -6.1096

$detailed.Range("B25").Value = -2.83936

$detailed.Range("B26").Value = -5.50985

$detailed.Range("B27").Value = -5.01

$detailed.Range("B28").Value = -0.89316

$detailed.Range("B29").Value = -5.01

$detailed.Range("B30").Value = 0.00002

$detailed.Range("B32").Value = 0.00848

$detailed.Range("B34").Value = 1.68233

$detailed.Range("B35").Value = -2.49239

$detailed.Range("B37").Value = 0.85381

$detailed.Range("B38").Value = 12.09485

$detailed.Range("B39").Value = 42.36398

$detailed.Range("B40").Value = 56.76084

$detailed.Range("B41").Value = 64.02197

$detailed.Range("B42").Value = 59.67565

$detailed.Range("B44").Value = 65

$detailed.Range("B45").Value = 62.50682

$detailed.Range("B46").Value = 61.52231

$detailed.Range("B47").Value = 58.9823

$detailed.Range("B48").Value = 58.49509

$detailed.Range("B49").Value = 62.01673
